$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: key = destination row number, value = source row number (from the
# pre-edit state). The edit is a permutation of data rows 4-18: the full
# per-species record (id, taxon order, redlist status, taxon id, names,
# author, coordinates, observers) moves to a different row position while
# the rest of each physical row (location, dates, etc.) stays put.
$mapping = @{
    4  = 16
    5  = 15
    6  = 17
    7  = 11
    8  = 7
    9  = 18
    10 = 14
    11 = 12
    12 = 10
    13 = 9
    14 = 13
    15 = 4
    16 = 6
    17 = 8
    18 = 5
}

$firstRow = 4
$lastRow = 18

# Only these columns actually move between rows; touching any other column
# (e.g. date-like text in Y/AA) risks Excel reinterpreting/reformatting it.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "AX")

# Snapshot every source row's values for the moved columns before any writes,
# so overlapping source/destination rows don't clobber data we still need.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $rowVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $rowVals[$col]
    }
}
